# The sheet lists (image_id, numEggs[, initials]) rows that should be kept
# sorted ascending by column A (image_id). A couple of rows had gotten
# shuffled out of order (a recycled/out-of-sequence "nadom#" id), so re-sort
# the whole data range to restore ascending order by column A - this is
# exactly what Data > Sort does in the Excel UI, and is what produced the
# <sortState>/<sortCondition> now persisted on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A2:C275")
$keyRange  = $ws.Range("A2:A275")

[void]$ws.Sort.SortFields.Clear()
[void]$ws.Sort.SortFields.Add($keyRange, 0, 1, $null, 0)
[void]$ws.Sort.SetRange($dataRange)
$ws.Sort.Header = 0
$ws.Sort.MatchCase = $false
$ws.Sort.Orientation = 1
[void]$ws.Sort.Apply()

# Leave the UI where the author ended up after fixing the id: the whole
# row for the corrected record selected, scrolled near it.
[void]$ws.Rows.Item(124).Select()
